$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct rows whose match data was reshuffled/fixed ("id" col A and "Date" col E are unchanged) ---

$row12a = New-Object "object[,]" 1,3
$row12a[0,0] = 6533424
$row12a[0,1] = 'Denmark Division 2'
$row12a[0,2] = 'Denmark Division 2'
$ws.Range("B12:D12").Value = $row12a
$row12b = New-Object "object[,]" 1,24
$row12b[0,0] = 'Kolding IF'
$row12b[0,1] = 'Esbjerg'
$row12b[0,2] = 1
$row12b[0,3] = 2
$row12b[0,4] = 'A'
$row12b[0,5] = 2.5
$row12b[0,6] = 3.75
$row12b[0,7] = 2.3
$row12b[0,8] = 2.25
$row12b[0,9] = 4
$row12b[0,10] = 2.5
$row12b[0,11] = 0
$row12b[0,12] = 1.775
$row12b[0,13] = 2.025
$row12b[0,14] = 2.75
$row12b[0,15] = 1.9
$row12b[0,16] = 1.9
$row12b[0,17] = -1
$row12b[0,18] = -1
$row12b[0,19] = 1.5
$row12b[0,20] = -1
$row12b[0,21] = 1.025
$row12b[0,22] = 0.45
$row12b[0,23] = -0.5
$ws.Range("F12:AC12").Value = $row12b

$row13a = New-Object "object[,]" 1,3
$row13a[0,0] = 6532919
$row13a[0,1] = 'Denmark Division 2'
$row13a[0,2] = 'Denmark Division 2'
$ws.Range("B13:D13").Value = $row13a
$row13b = New-Object "object[,]" 1,24
$row13b[0,0] = 'Aarhus Fremad'
$row13b[0,1] = 'AB Copenhagen'
$row13b[0,2] = 5
$row13b[0,3] = 2
$row13b[0,4] = 'H'
$row13b[0,5] = 1.727
$row13b[0,6] = 3.8
$row13b[0,7] = 3.75
$row13b[0,8] = 1.5
$row13b[0,9] = 4.333
$row13b[0,10] = 5
$row13b[0,11] = -1
$row13b[0,12] = 1.8
$row13b[0,13] = 2
$row13b[0,14] = 3.25
$row13b[0,15] = 1.95
$row13b[0,16] = 1.85
$row13b[0,17] = 0.5
$row13b[0,18] = -1
$row13b[0,19] = -1
$row13b[0,20] = 0.8
$row13b[0,21] = -1
$row13b[0,22] = 0.95
$row13b[0,23] = -1
$ws.Range("F13:AC13").Value = $row13b

$row14a = New-Object "object[,]" 1,3
$row14a[0,0] = 6529391
$row14a[0,1] = 'Denmark Division 2'
$row14a[0,2] = 'Denmark Division 2'
$ws.Range("B14:D14").Value = $row14a
$row14b = New-Object "object[,]" 1,24
$row14b[0,0] = 'B93 Copenhagen'
$row14b[0,1] = 'Thisted FC'
$row14b[0,2] = 1
$row14b[0,3] = 1
$row14b[0,4] = 'D'
$row14b[0,5] = 1.3
$row14b[0,6] = 4.5
$row14b[0,7] = 9
$row14b[0,8] = 1.285
$row14b[0,9] = 5.25
$row14b[0,10] = 9
$row14b[0,11] = -1.5
$row14b[0,12] = 1.825
$row14b[0,13] = 1.975
$row14b[0,14] = 3
$row14b[0,15] = 1.925
$row14b[0,16] = 1.875
$row14b[0,17] = -1
$row14b[0,18] = 4.25
$row14b[0,19] = -1
$row14b[0,20] = -1
$row14b[0,21] = 0.9750000000000001
$row14b[0,22] = -1
$row14b[0,23] = 0.875
$ws.Range("F14:AC14").Value = $row14b

$row15a = New-Object "object[,]" 1,3
$row15a[0,0] = 6529284
$row15a[0,1] = 'Denmark Division 2'
$row15a[0,2] = 'Denmark Division 2'
$ws.Range("B15:D15").Value = $row15a
$row15b = New-Object "object[,]" 1,24
$row15b[0,0] = 'Frem'
$row15b[0,1] = 'Brabrand'
$row15b[0,2] = 1
$row15b[0,3] = 1
$row15b[0,4] = 'D'
$row15b[0,5] = 2.375
$row15b[0,6] = 3.8
$row15b[0,7] = 2.375
$row15b[0,8] = 2.4
$row15b[0,9] = 3.8
$row15b[0,10] = 2.3
$row15b[0,11] = 0
$row15b[0,12] = 1.975
$row15b[0,13] = 1.825
$row15b[0,14] = 3
$row15b[0,15] = 1.95
$row15b[0,16] = 1.85
$row15b[0,17] = -1
$row15b[0,18] = 2.8
$row15b[0,19] = -1
$row15b[0,20] = 0
$row15b[0,21] = 0
$row15b[0,22] = -1
$row15b[0,23] = 0.8500000000000001
$ws.Range("F15:AC15").Value = $row15b

$row17a = New-Object "object[,]" 1,3
$row17a[0,0] = 6858901
$row17a[0,1] = 'Denmark Division 2'
$row17a[0,2] = 'Denmark Division 2'
$ws.Range("B17:D17").Value = $row17a
$row17b = New-Object "object[,]" 1,24
$row17b[0,0] = 'Skive'
$row17b[0,1] = 'Esbjerg'
$row17b[0,2] = 1
$row17b[0,3] = 4
$row17b[0,4] = 'A'
$row17b[0,5] = 3.75
$row17b[0,6] = 3.75
$row17b[0,7] = 1.727
$row17b[0,8] = 4.333
$row17b[0,9] = 4
$row17b[0,10] = 1.6
$row17b[0,11] = 0.75
$row17b[0,12] = 2
$row17b[0,13] = 1.8
$row17b[0,14] = 3
$row17b[0,15] = 1.95
$row17b[0,16] = 1.85
$row17b[0,17] = -1
$row17b[0,18] = -1
$row17b[0,19] = 0.6000000000000001
$row17b[0,20] = -1
$row17b[0,21] = 0.8
$row17b[0,22] = 0.95
$row17b[0,23] = -1
$ws.Range("F17:AC17").Value = $row17b

$row18a = New-Object "object[,]" 1,3
$row18a[0,0] = 6858900
$row18a[0,1] = 'Denmark Division 2'
$row18a[0,2] = 'Denmark Division 2'
$ws.Range("B18:D18").Value = $row18a
$row18b = New-Object "object[,]" 1,24
$row18b[0,0] = 'Fremad Amager'
$row18b[0,1] = 'Hellerup IK'
$row18b[0,2] = 2
$row18b[0,3] = 1
$row18b[0,4] = 'H'
$row18b[0,5] = 1.833
$row18b[0,6] = 3.6
$row18b[0,7] = 3.4
$row18b[0,8] = 1.833
$row18b[0,9] = 3.6
$row18b[0,10] = 3.4
$row18b[0,11] = -0.5
$row18b[0,12] = 1.875
$row18b[0,13] = 1.925
$row18b[0,14] = 2.75
$row18b[0,15] = 1.825
$row18b[0,16] = 1.975
$row18b[0,17] = 0.833
$row18b[0,18] = -1
$row18b[0,19] = -1
$row18b[0,20] = 0.875
$row18b[0,21] = -1
$row18b[0,22] = 0.4125
$row18b[0,23] = -0.5
$ws.Range("F18:AC18").Value = $row18b

$row29a = New-Object "object[,]" 1,3
$row29a[0,0] = 6858913
$row29a[0,1] = 'Denmark Division 2'
$row29a[0,2] = 'Denmark Division 2'
$ws.Range("B29:D29").Value = $row29a
$row29b = New-Object "object[,]" 1,24
$row29b[0,0] = 'Nykobing'
$row29b[0,1] = 'FC Roskilde'
$row29b[0,2] = 1
$row29b[0,3] = 2
$row29b[0,4] = 'A'
$row29b[0,5] = 2.05
$row29b[0,6] = 3.5
$row29b[0,7] = 3.1
$row29b[0,8] = 2.15
$row29b[0,9] = 3.4
$row29b[0,10] = 2.9
$row29b[0,11] = -0.25
$row29b[0,12] = 1.925
$row29b[0,13] = 1.875
$row29b[0,14] = 2.75
$row29b[0,15] = 1.875
$row29b[0,16] = 1.925
$row29b[0,17] = -1
$row29b[0,18] = -1
$row29b[0,19] = 1.9
$row29b[0,20] = -1
$row29b[0,21] = 0.875
$row29b[0,22] = 0.4375
$row29b[0,23] = -0.5
$ws.Range("F29:AC29").Value = $row29b

$row30a = New-Object "object[,]" 1,3
$row30a[0,0] = 6858912
$row30a[0,1] = 'Denmark Division 2'
$row30a[0,2] = 'Denmark Division 2'
$ws.Range("B30:D30").Value = $row30a
$row30b = New-Object "object[,]" 1,24
$row30b[0,0] = 'AB Copenhagen'
$row30b[0,1] = 'Hellerup IK'
$row30b[0,2] = 5
$row30b[0,3] = 4
$row30b[0,4] = 'H'
$row30b[0,5] = 1.85
$row30b[0,6] = 3.5
$row30b[0,7] = 3.75
$row30b[0,8] = 1.85
$row30b[0,9] = 3.5
$row30b[0,10] = 3.75
$row30b[0,11] = -0.5
$row30b[0,12] = 1.875
$row30b[0,13] = 1.925
$row30b[0,14] = 3
$row30b[0,15] = 1.925
$row30b[0,16] = 1.875
$row30b[0,17] = 0.8500000000000001
$row30b[0,18] = -1
$row30b[0,19] = -1
$row30b[0,20] = 0.875
$row30b[0,21] = -1
$row30b[0,22] = 0.925
$row30b[0,23] = -1
$ws.Range("F30:AC30").Value = $row30b

$row31a = New-Object "object[,]" 1,3
$row31a[0,0] = 6858911
$row31a[0,1] = 'Denmark Division 2'
$row31a[0,2] = 'Denmark Division 2'
$ws.Range("B31:D31").Value = $row31a
$row31b = New-Object "object[,]" 1,24
$row31b[0,0] = 'Skive'
$row31b[0,1] = 'Aarhus Fremad'
$row31b[0,2] = 1
$row31b[0,3] = 2
$row31b[0,4] = 'A'
$row31b[0,5] = 4.75
$row31b[0,6] = 3.8
$row31b[0,7] = 1.6
$row31b[0,8] = 5
$row31b[0,9] = 3.8
$row31b[0,10] = 1.571
$row31b[0,11] = 1
$row31b[0,12] = 1.8
$row31b[0,13] = 2
$row31b[0,14] = 2.75
$row31b[0,15] = 1.75
$row31b[0,16] = 1.95
$row31b[0,17] = -1
$row31b[0,18] = -1
$row31b[0,19] = 0.571
$row31b[0,20] = 0
$row31b[0,21] = 0
$row31b[0,22] = 0.375
$row31b[0,23] = -0.5
$ws.Range("F31:AC31").Value = $row31b

$row46a = New-Object "object[,]" 1,3
$row46a[0,0] = 6858930
$row46a[0,1] = 'Denmark Division 2'
$row46a[0,2] = 'Denmark Division 2'
$ws.Range("B46:D46").Value = $row46a
$row46b = New-Object "object[,]" 1,24
$row46b[0,0] = 'Thisted FC'
$row46b[0,1] = 'Skive'
$row46b[0,2] = 0
$row46b[0,3] = 0
$row46b[0,4] = 'D'
$row46b[0,5] = 2.4
$row46b[0,6] = 3.5
$row46b[0,7] = 2.4
$row46b[0,8] = 2.25
$row46b[0,9] = 3.5
$row46b[0,10] = 2.55
$row46b[0,11] = 0
$row46b[0,12] = 1.775
$row46b[0,13] = 2.025
$row46b[0,14] = 2.5
$row46b[0,15] = 1.85
$row46b[0,16] = 1.95
$row46b[0,17] = -1
$row46b[0,18] = 2.5
$row46b[0,19] = -1
$row46b[0,20] = 0
$row46b[0,21] = 0
$row46b[0,22] = -1
$row46b[0,23] = 0.95
$ws.Range("F46:AC46").Value = $row46b

$row47a = New-Object "object[,]" 1,3
$row47a[0,0] = 6858929
$row47a[0,1] = 'Denmark Division 2'
$row47a[0,2] = 'Denmark Division 2'
$ws.Range("B47:D47").Value = $row47a
$row47b = New-Object "object[,]" 1,24
$row47b[0,0] = 'FC Roskilde'
$row47b[0,1] = 'AB Copenhagen'
$row47b[0,2] = 2
$row47b[0,3] = 0
$row47b[0,4] = 'H'
$row47b[0,5] = 2.45
$row47b[0,6] = 3.5
$row47b[0,7] = 2.45
$row47b[0,8] = 2.375
$row47b[0,9] = 3.4
$row47b[0,10] = 2.55
$row47b[0,11] = 0
$row47b[0,12] = 1.85
$row47b[0,13] = 1.95
$row47b[0,14] = 2.5
$row47b[0,15] = 1.875
$row47b[0,16] = 1.925
$row47b[0,17] = 1.375
$row47b[0,18] = -1
$row47b[0,19] = -1
$row47b[0,20] = 0.8500000000000001
$row47b[0,21] = -1
$row47b[0,22] = -1
$row47b[0,23] = 0.925
$ws.Range("F47:AC47").Value = $row47b

$row48a = New-Object "object[,]" 1,3
$row48a[0,0] = 6858928
$row48a[0,1] = 'Denmark Division 2'
$row48a[0,2] = 'Denmark Division 2'
$ws.Range("B48:D48").Value = $row48a
$row48b = New-Object "object[,]" 1,24
$row48b[0,0] = 'FA 2000'
$row48b[0,1] = 'Hellerup IK'
$row48b[0,2] = 2
$row48b[0,3] = 0
$row48b[0,4] = 'H'
$row48b[0,5] = 2.2
$row48b[0,6] = 3.5
$row48b[0,7] = 2.8
$row48b[0,8] = 2.1
$row48b[0,9] = 3.4
$row48b[0,10] = 3.1
$row48b[0,11] = -0.25
$row48b[0,12] = 1.85
$row48b[0,13] = 1.95
$row48b[0,14] = 2.75
$row48b[0,15] = 1.925
$row48b[0,16] = 1.875
$row48b[0,17] = 1.1
$row48b[0,18] = -1
$row48b[0,19] = -1
$row48b[0,20] = 0.8500000000000001
$row48b[0,21] = -1
$row48b[0,22] = -1
$row48b[0,23] = 0.875
$ws.Range("F48:AC48").Value = $row48b

$row66a = New-Object "object[,]" 1,3
$row66a[0,0] = 6858950
$row66a[0,1] = 'Denmark Division 2'
$row66a[0,2] = 'Denmark Division 2'
$ws.Range("B66:D66").Value = $row66a
$row66b = New-Object "object[,]" 1,24
$row66b[0,0] = 'FC Roskilde'
$row66b[0,1] = 'Esbjerg'
$row66b[0,2] = 3
$row66b[0,3] = 3
$row66b[0,4] = 'D'
$row66b[0,5] = 4
$row66b[0,6] = 3.5
$row66b[0,7] = 1.8
$row66b[0,8] = 4
$row66b[0,9] = 3.5
$row66b[0,10] = 1.8
$row66b[0,11] = 0.5
$row66b[0,12] = 1.975
$row66b[0,13] = 1.825
$row66b[0,14] = 2.75
$row66b[0,15] = 1.825
$row66b[0,16] = 1.975
$row66b[0,17] = -1
$row66b[0,18] = 2.5
$row66b[0,19] = -1
$row66b[0,20] = 0.9750000000000001
$row66b[0,21] = -1
$row66b[0,22] = 0.825
$row66b[0,23] = -1
$ws.Range("F66:AC66").Value = $row66b

$row67a = New-Object "object[,]" 1,3
$row67a[0,0] = 6858951
$row67a[0,1] = 'Denmark Division 2'
$row67a[0,2] = 'Denmark Division 2'
$ws.Range("B67:D67").Value = $row67a
$row67b = New-Object "object[,]" 1,24
$row67b[0,0] = 'Brabrand'
$row67b[0,1] = 'Skive'
$row67b[0,2] = 1
$row67b[0,3] = 0
$row67b[0,4] = 'H'
$row67b[0,5] = 2.55
$row67b[0,6] = 3.5
$row67b[0,7] = 2.3
$row67b[0,8] = 2.55
$row67b[0,9] = 3.4
$row67b[0,10] = 2.375
$row67b[0,11] = 0
$row67b[0,12] = 1.975
$row67b[0,13] = 1.825
$row67b[0,14] = 2.25
$row67b[0,15] = 1.825
$row67b[0,16] = 1.975
$row67b[0,17] = 1.55
$row67b[0,18] = -1
$row67b[0,19] = -1
$row67b[0,20] = 0.9750000000000001
$row67b[0,21] = -1
$row67b[0,22] = -1
$row67b[0,23] = 0.9750000000000001
$ws.Range("F67:AC67").Value = $row67b

$row71a = New-Object "object[,]" 1,3
$row71a[0,0] = 6858955
$row71a[0,1] = 'Denmark Division 2'
$row71a[0,2] = 'Denmark Division 2'
$ws.Range("B71:D71").Value = $row71a
$row71b = New-Object "object[,]" 1,24
$row71b[0,0] = 'Aarhus Fremad'
$row71b[0,1] = 'AB Copenhagen'
$row71b[0,2] = 1
$row71b[0,3] = 1
$row71b[0,4] = 'D'
$row71b[0,5] = 1.6
$row71b[0,6] = 3.75
$row71b[0,7] = 4.75
$row71b[0,8] = 1.444
$row71b[0,9] = 4
$row71b[0,10] = 6.5
$row71b[0,11] = -1.25
$row71b[0,12] = 1.975
$row71b[0,13] = 1.825
$row71b[0,14] = 3
$row71b[0,15] = 1.8
$row71b[0,16] = 2
$row71b[0,17] = -1
$row71b[0,18] = 3
$row71b[0,19] = -1
$row71b[0,20] = -1
$row71b[0,21] = 0.825
$row71b[0,22] = -1
$row71b[0,23] = 1
$ws.Range("F71:AC71").Value = $row71b

$row72a = New-Object "object[,]" 1,3
$row72a[0,0] = 6858954
$row72a[0,1] = 'Denmark Division 2'
$row72a[0,2] = 'Denmark Division 2'
$ws.Range("B72:D72").Value = $row72a
$row72b = New-Object "object[,]" 1,24
$row72b[0,0] = 'Middelfart'
$row72b[0,1] = 'FC Roskilde'
$row72b[0,2] = 0
$row72b[0,3] = 1
$row72b[0,4] = 'A'
$row72b[0,5] = 2.15
$row72b[0,6] = 3.5
$row72b[0,7] = 2.875
$row72b[0,8] = 2.2
$row72b[0,9] = 3.4
$row72b[0,10] = 2.875
$row72b[0,11] = -0.25
$row72b[0,12] = 1.975
$row72b[0,13] = 1.825
$row72b[0,14] = 2.5
$row72b[0,15] = 1.975
$row72b[0,16] = 1.825
$row72b[0,17] = -1
$row72b[0,18] = -1
$row72b[0,19] = 1.875
$row72b[0,20] = -1
$row72b[0,21] = 0.825
$row72b[0,22] = -1
$row72b[0,23] = 0.825
$ws.Range("F72:AC72").Value = $row72b

$row82a = New-Object "object[,]" 1,3
$row82a[0,0] = 6859010
$row82a[0,1] = 'Denmark Division 2'
$row82a[0,2] = 'Denmark Division 2'
$ws.Range("B82:D82").Value = $row82a
$row82b = New-Object "object[,]" 1,24
$row82b[0,0] = 'Esbjerg'
$row82b[0,1] = 'FA 2000'
$row82b[0,2] = 3
$row82b[0,3] = 1
$row82b[0,4] = 'H'
$row82b[0,5] = 1.222
$row82b[0,6] = 6.5
$row82b[0,7] = 9
$row82b[0,8] = 1.2
$row82b[0,9] = 6.5
$row82b[0,10] = 10
$row82b[0,11] = -2
$row82b[0,12] = 1.9
$row82b[0,13] = 1.9
$row82b[0,14] = 3.5
$row82b[0,15] = 1.875
$row82b[0,16] = 1.925
$row82b[0,17] = 0.2
$row82b[0,18] = -1
$row82b[0,19] = -1
$row82b[0,20] = 0
$row82b[0,21] = 0
$row82b[0,22] = 0.875
$row82b[0,23] = -1
$ws.Range("F82:AC82").Value = $row82b

$row84a = New-Object "object[,]" 1,3
$row84a[0,0] = 6859007
$row84a[0,1] = 'Denmark Division 2'
$row84a[0,2] = 'Denmark Division 2'
$ws.Range("B84:D84").Value = $row84a
$row84b = New-Object "object[,]" 1,24
$row84b[0,0] = 'Skive'
$row84b[0,1] = 'FC Roskilde'
$row84b[0,2] = 1
$row84b[0,3] = 2
$row84b[0,4] = 'A'
$row84b[0,5] = 3.6
$row84b[0,6] = 3.4
$row84b[0,7] = 1.909
$row84b[0,8] = 3.2
$row84b[0,9] = 3.4
$row84b[0,10] = 2.05
$row84b[0,11] = 0.25
$row84b[0,12] = 2
$row84b[0,13] = 1.8
$row84b[0,14] = 2.75
$row84b[0,15] = 1.975
$row84b[0,16] = 1.825
$row84b[0,17] = -1
$row84b[0,18] = -1
$row84b[0,19] = 1.05
$row84b[0,20] = -1
$row84b[0,21] = 0.8
$row84b[0,22] = 0.4875
$row84b[0,23] = -0.5
$ws.Range("F84:AC84").Value = $row84b

$row85a = New-Object "object[,]" 1,3
$row85a[0,0] = 6859011
$row85a[0,1] = 'Denmark Division 2'
$row85a[0,2] = 'Denmark Division 2'
$ws.Range("B85:D85").Value = $row85a
$row85b = New-Object "object[,]" 1,24
$row85b[0,0] = 'Middelfart'
$row85b[0,1] = 'Nykobing'
$row85b[0,2] = 2
$row85b[0,3] = 2
$row85b[0,4] = 'D'
$row85b[0,5] = 2
$row85b[0,6] = 3.5
$row85b[0,7] = 3.3
$row85b[0,8] = 2
$row85b[0,9] = 3.5
$row85b[0,10] = 3.3
$row85b[0,11] = -0.25
$row85b[0,12] = 1.75
$row85b[0,13] = 1.95
$row85b[0,14] = 2.5
$row85b[0,15] = 1.8
$row85b[0,16] = 2
$row85b[0,17] = -1
$row85b[0,18] = 2.5
$row85b[0,19] = -1
$row85b[0,20] = -0.5
$row85b[0,21] = 0.475
$row85b[0,22] = 0.8
$row85b[0,23] = -1
$ws.Range("F85:AC85").Value = $row85b

$row92a = New-Object "object[,]" 1,3
$row92a[0,0] = 6859028
$row92a[0,1] = 'Denmark Division 2'
$row92a[0,2] = 'Denmark Division 2'
$ws.Range("B92:D92").Value = $row92a
$row92b = New-Object "object[,]" 1,24
$row92b[0,0] = 'Thisted FC'
$row92b[0,1] = 'AB Copenhagen'
$row92b[0,2] = 0
$row92b[0,3] = 1
$row92b[0,4] = 'A'
$row92b[0,5] = 2.5
$row92b[0,6] = 3.4
$row92b[0,7] = 2.5
$row92b[0,8] = 2.5
$row92b[0,9] = 3.4
$row92b[0,10] = 2.5
$row92b[0,11] = 0
$row92b[0,12] = 1.9
$row92b[0,13] = 1.9
$row92b[0,14] = 2.75
$row92b[0,15] = 1.9
$row92b[0,16] = 1.9
$row92b[0,17] = -1
$row92b[0,18] = -1
$row92b[0,19] = 1.5
$row92b[0,20] = -1
$row92b[0,21] = 0.8999999999999999
$row92b[0,22] = -1
$row92b[0,23] = 0.8999999999999999
$ws.Range("F92:AC92").Value = $row92b

$row93a = New-Object "object[,]" 1,3
$row93a[0,0] = 6859026
$row93a[0,1] = 'Denmark Division 2'
$row93a[0,2] = 'Denmark Division 2'
$ws.Range("B93:D93").Value = $row93a
$row93b = New-Object "object[,]" 1,24
$row93b[0,0] = 'Nykobing'
$row93b[0,1] = 'Esbjerg'
$row93b[0,2] = 0
$row93b[0,3] = 3
$row93b[0,4] = 'A'
$row93b[0,5] = 4.333
$row93b[0,6] = 4
$row93b[0,7] = 1.615
$row93b[0,8] = 4.5
$row93b[0,9] = 4
$row93b[0,10] = 1.571
$row93b[0,11] = 1
$row93b[0,12] = 1.825
$row93b[0,13] = 1.975
$row93b[0,14] = 3.25
$row93b[0,15] = 1.9
$row93b[0,16] = 1.9
$row93b[0,17] = -1
$row93b[0,18] = -1
$row93b[0,19] = 0.571
$row93b[0,20] = -1
$row93b[0,21] = 0.9750000000000001
$row93b[0,22] = -0.5
$row93b[0,23] = 0.45
$ws.Range("F93:AC93").Value = $row93b

$row95a = New-Object "object[,]" 1,3
$row95a[0,0] = 6859041
$row95a[0,1] = 'Denmark Division 2'
$row95a[0,2] = 'Denmark Division 2'
$ws.Range("B95:D95").Value = $row95a
$row95b = New-Object "object[,]" 1,24
$row95b[0,0] = 'Middelfart'
$row95b[0,1] = 'FA 2000'
$row95b[0,2] = 2
$row95b[0,3] = 1
$row95b[0,4] = 'H'
$row95b[0,5] = 1.65
$row95b[0,6] = 3.6
$row95b[0,7] = 4.5
$row95b[0,8] = 1.571
$row95b[0,9] = 3.75
$row95b[0,10] = 5.25
$row95b[0,11] = -0.75
$row95b[0,12] = 1.775
$row95b[0,13] = 2.025
$row95b[0,14] = 2.5
$row95b[0,15] = 1.9
$row95b[0,16] = 1.9
$row95b[0,17] = 0.571
$row95b[0,18] = -1
$row95b[0,19] = -1
$row95b[0,20] = 0.3875
$row95b[0,21] = -0.5
$row95b[0,22] = 0.8999999999999999
$row95b[0,23] = -1
$ws.Range("F95:AC95").Value = $row95b

$row97a = New-Object "object[,]" 1,3
$row97a[0,0] = 6859037
$row97a[0,1] = 'Denmark Division 2'
$row97a[0,2] = 'Denmark Division 2'
$ws.Range("B97:D97").Value = $row97a
$row97b = New-Object "object[,]" 1,24
$row97b[0,0] = 'Aarhus Fremad'
$row97b[0,1] = 'Hellerup IK'
$row97b[0,2] = 3
$row97b[0,3] = 2
$row97b[0,4] = 'H'
$row97b[0,5] = 1.5
$row97b[0,6] = 4
$row97b[0,7] = 5.5
$row97b[0,8] = 1.363
$row97b[0,9] = 4.333
$row97b[0,10] = 7.5
$row97b[0,11] = -1.25
$row97b[0,12] = 1.775
$row97b[0,13] = 2.025
$row97b[0,14] = 3.25
$row97b[0,15] = 1.925
$row97b[0,16] = 1.875
$row97b[0,17] = 0.363
$row97b[0,18] = -1
$row97b[0,19] = -1
$row97b[0,20] = -0.5
$row97b[0,21] = 0.5125
$row97b[0,22] = 0.925
$row97b[0,23] = -1
$ws.Range("F97:AC97").Value = $row97b

$row107a = New-Object "object[,]" 1,3
$row107a[0,0] = 6859067
$row107a[0,1] = 'Denmark Division 2'
$row107a[0,2] = 'Denmark Division 2'
$ws.Range("B107:D107").Value = $row107a
$row107b = New-Object "object[,]" 1,24
$row107b[0,0] = 'AB Copenhagen'
$row107b[0,1] = 'FA 2000'
$row107b[0,2] = 3
$row107b[0,3] = 2
$row107b[0,4] = 'H'
$row107b[0,5] = 1.727
$row107b[0,6] = 3.6
$row107b[0,7] = 4.2
$row107b[0,8] = 1.7
$row107b[0,9] = 3.6
$row107b[0,10] = 4.5
$row107b[0,11] = -0.75
$row107b[0,12] = 1.85
$row107b[0,13] = 1.95
$row107b[0,14] = 2.75
$row107b[0,15] = 1.75
$row107b[0,16] = 1.95
$row107b[0,17] = 0.7
$row107b[0,18] = -1
$row107b[0,19] = -1
$row107b[0,20] = 0.425
$row107b[0,21] = -0.5
$row107b[0,22] = 0.75
$row107b[0,23] = -1
$ws.Range("F107:AC107").Value = $row107b

$row108a = New-Object "object[,]" 1,3
$row108a[0,0] = 6859059
$row108a[0,1] = 'Denmark Division 2'
$row108a[0,2] = 'Denmark Division 2'
$ws.Range("B108:D108").Value = $row108a
$row108b = New-Object "object[,]" 1,24
$row108b[0,0] = 'Hellerup IK'
$row108b[0,1] = 'Esbjerg'
$row108b[0,2] = 1
$row108b[0,3] = 4
$row108b[0,4] = 'A'
$row108b[0,5] = 7
$row108b[0,6] = 5
$row108b[0,7] = 1.333
$row108b[0,8] = 8
$row108b[0,9] = 5.25
$row108b[0,10] = 1.3
$row108b[0,11] = 1.5
$row108b[0,12] = 1.925
$row108b[0,13] = 1.875
$row108b[0,14] = 3.25
$row108b[0,15] = 1.875
$row108b[0,16] = 1.925
$row108b[0,17] = -1
$row108b[0,18] = -1
$row108b[0,19] = 0.3
$row108b[0,20] = -1
$row108b[0,21] = 0.875
$row108b[0,22] = 0.875
$row108b[0,23] = -1
$ws.Range("F108:AC108").Value = $row108b

$row113a = New-Object "object[,]" 1,3
$row113a[0,0] = 6859074
$row113a[0,1] = 'Denmark Division 2'
$row113a[0,2] = 'Denmark Division 2'
$ws.Range("B113:D113").Value = $row113a
$row113b = New-Object "object[,]" 1,24
$row113b[0,0] = 'FC Roskilde'
$row113b[0,1] = 'Middelfart'
$row113b[0,2] = 3
$row113b[0,3] = 1
$row113b[0,4] = 'H'
$row113b[0,5] = 2.2
$row113b[0,6] = 3.4
$row113b[0,7] = 2.8
$row113b[0,8] = 2.15
$row113b[0,9] = 3.4
$row113b[0,10] = 2.875
$row113b[0,11] = -0.25
$row113b[0,12] = 1.95
$row113b[0,13] = 1.85
$row113b[0,14] = 2.5
$row113b[0,15] = 1.975
$row113b[0,16] = 1.825
$row113b[0,17] = 1.15
$row113b[0,18] = -1
$row113b[0,19] = -1
$row113b[0,20] = 0.95
$row113b[0,21] = -1
$row113b[0,22] = 0.9750000000000001
$row113b[0,23] = -1
$ws.Range("F113:AC113").Value = $row113b

$row114a = New-Object "object[,]" 1,3
$row114a[0,0] = 6859073
$row114a[0,1] = 'Denmark Division 2'
$row114a[0,2] = 'Denmark Division 2'
$ws.Range("B114:D114").Value = $row114a
$row114b = New-Object "object[,]" 1,24
$row114b[0,0] = 'Aarhus Fremad'
$row114b[0,1] = 'Thisted FC'
$row114b[0,2] = 1
$row114b[0,3] = 0
$row114b[0,4] = 'H'
$row114b[0,5] = 1.45
$row114b[0,6] = 4
$row114b[0,7] = 6
$row114b[0,8] = 1.333
$row114b[0,9] = 4.75
$row114b[0,10] = 7.5
$row114b[0,11] = -1.5
$row114b[0,12] = 1.95
$row114b[0,13] = 1.85
$row114b[0,14] = 3
$row114b[0,15] = 1.925
$row114b[0,16] = 1.875
$row114b[0,17] = 0.333
$row114b[0,18] = -1
$row114b[0,19] = -1
$row114b[0,20] = -1
$row114b[0,21] = 0.8500000000000001
$row114b[0,22] = -1
$row114b[0,23] = 0.875
$ws.Range("F114:AC114").Value = $row114b

$row118a = New-Object "object[,]" 1,3
$row118a[0,0] = 6859066
$row118a[0,1] = 'Denmark Division 2'
$row118a[0,2] = 'Denmark Division 2'
$ws.Range("B118:D118").Value = $row118a
$row118b = New-Object "object[,]" 1,24
$row118b[0,0] = 'Hellerup IK'
$row118b[0,1] = 'FC Roskilde'
$row118b[0,2] = 2
$row118b[0,3] = 3
$row118b[0,4] = 'A'
$row118b[0,5] = 3.6
$row118b[0,6] = 3.6
$row118b[0,7] = 1.833
$row118b[0,8] = 4.5
$row118b[0,9] = 3.8
$row118b[0,10] = 1.615
$row118b[0,11] = 0.75
$row118b[0,12] = 1.975
$row118b[0,13] = 1.825
$row118b[0,14] = 2.75
$row118b[0,15] = 1.875
$row118b[0,16] = 1.925
$row118b[0,17] = -1
$row118b[0,18] = -1
$row118b[0,19] = 0.615
$row118b[0,20] = -0.5
$row118b[0,21] = 0.4125
$row118b[0,22] = 0.875
$row118b[0,23] = -1
$ws.Range("F118:AC118").Value = $row118b

$row119a = New-Object "object[,]" 1,3
$row119a[0,0] = 6859068
$row119a[0,1] = 'Denmark Division 2'
$row119a[0,2] = 'Denmark Division 2'
$ws.Range("B119:D119").Value = $row119a
$row119b = New-Object "object[,]" 1,24
$row119b[0,0] = 'Fremad Amager'
$row119b[0,1] = 'FA 2000'
$row119b[0,2] = 0
$row119b[0,3] = 0
$row119b[0,4] = 'D'
$row119b[0,5] = 2.2
$row119b[0,6] = 3.4
$row119b[0,7] = 2.9
$row119b[0,8] = 2.15
$row119b[0,9] = 3.4
$row119b[0,10] = 2.9
$row119b[0,11] = -0.25
$row119b[0,12] = 1.925
$row119b[0,13] = 1.875
$row119b[0,14] = 2.5
$row119b[0,15] = 1.925
$row119b[0,16] = 1.875
$row119b[0,17] = -1
$row119b[0,18] = 2.4
$row119b[0,19] = -1
$row119b[0,20] = -0.5
$row119b[0,21] = 0.4375
$row119b[0,22] = -1
$row119b[0,23] = 0.875
$ws.Range("F119:AC119").Value = $row119b

$row132a = New-Object "object[,]" 1,3
$row132a[0,0] = 6859042
$row132a[0,1] = 'Denmark Division 2'
$row132a[0,2] = 'Denmark Division 2'
$ws.Range("B132:D132").Value = $row132a
$row132b = New-Object "object[,]" 1,24
$row132b[0,0] = 'Aarhus Fremad'
$row132b[0,1] = 'Skive'
$row132b[0,2] = 0
$row132b[0,3] = 1
$row132b[0,4] = 'A'
$row132b[0,5] = 1.533
$row132b[0,6] = 3.8
$row132b[0,7] = 5
$row132b[0,8] = 1.65
$row132b[0,9] = 3.5
$row132b[0,10] = 4.333
$row132b[0,11] = -0.75
$row132b[0,12] = 1.9
$row132b[0,13] = 1.9
$row132b[0,14] = 2.5
$row132b[0,15] = 1.9
$row132b[0,16] = 1.9
$row132b[0,17] = -1
$row132b[0,18] = -1
$row132b[0,19] = 3.333
$row132b[0,20] = -1
$row132b[0,21] = 0.8999999999999999
$row132b[0,22] = -1
$row132b[0,23] = 0.8999999999999999
$ws.Range("F132:AC132").Value = $row132b

$row133a = New-Object "object[,]" 1,3
$row133a[0,0] = 6859040
$row133a[0,1] = 'Denmark Division 2'
$row133a[0,2] = 'Denmark Division 2'
$ws.Range("B133:D133").Value = $row133a
$row133b = New-Object "object[,]" 1,24
$row133b[0,0] = 'Esbjerg'
$row133b[0,1] = 'FC Roskilde'
$row133b[0,2] = 4
$row133b[0,3] = 2
$row133b[0,4] = 'H'
$row133b[0,5] = 1.65
$row133b[0,6] = 3.75
$row133b[0,7] = 4.333
$row133b[0,8] = 1.6
$row133b[0,9] = 3.8
$row133b[0,10] = 4.75
$row133b[0,11] = -1
$row133b[0,12] = 1.95
$row133b[0,13] = 1.75
$row133b[0,14] = 3
$row133b[0,15] = 1.825
$row133b[0,16] = 1.975
$row133b[0,17] = 0.6000000000000001
$row133b[0,18] = -1
$row133b[0,19] = -1
$row133b[0,20] = 0.95
$row133b[0,21] = -1
$row133b[0,22] = 0.825
$row133b[0,23] = -1
$ws.Range("F133:AC133").Value = $row133b

# --- Append brand-new rows 136-140 (copy row 135 formatting first, for the bordered id column and date format) ---
$ws.Range("A135:AC135").Copy()
$ws.Range("A136:AC136").PasteSpecial(-4122)
$ws.Range("A137:AC137").PasteSpecial(-4122)
$ws.Range("A138:AC138").PasteSpecial(-4122)
$ws.Range("A139:AC139").PasteSpecial(-4122)
$ws.Range("A140:AC140").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$row136 = New-Object "object[,]" 1,29
$row136[0,0] = 134
$row136[0,1] = 6859034
$row136[0,2] = 'Denmark Division 2'
$row136[0,3] = 'Denmark Division 2'
$row136[0,4] = 45383.33333333334
$row136[0,5] = 'Skive'
$row136[0,6] = 'Hellerup IK'
$row136[0,7] = $null
$row136[0,8] = $null
$row136[0,9] = $null
$row136[0,10] = 1.833
$row136[0,11] = 3.6
$row136[0,12] = 3.7
$row136[0,13] = 1.75
$row136[0,14] = 3.75
$row136[0,15] = 4
$row136[0,16] = -0.5
$row136[0,17] = 1.775
$row136[0,18] = 2.025
$row136[0,19] = 2.5
$row136[0,20] = 1.8
$row136[0,21] = 2
$row136[0,22] = 0
$row136[0,23] = 0
$row136[0,24] = 0
$row136[0,25] = 0
$row136[0,26] = 0
$row136[0,27] = $null
$row136[0,28] = $null
$ws.Range("A136:AC136").Value = $row136

$row137 = New-Object "object[,]" 1,29
$row137[0,0] = 135
$row137[0,1] = 6859030
$row137[0,2] = 'Denmark Division 2'
$row137[0,3] = 'Denmark Division 2'
$row137[0,4] = 45383.375
$row137[0,5] = 'FC Roskilde'
$row137[0,6] = 'Aarhus Fremad'
$row137[0,7] = $null
$row137[0,8] = $null
$row137[0,9] = $null
$row137[0,10] = 2.1
$row137[0,11] = 3.5
$row137[0,12] = 2.9
$row137[0,13] = 2.35
$row137[0,14] = 3.4
$row137[0,15] = 2.7
$row137[0,16] = 0
$row137[0,17] = 1.8
$row137[0,18] = 2
$row137[0,19] = 2.75
$row137[0,20] = 1.95
$row137[0,21] = 1.85
$row137[0,22] = 0
$row137[0,23] = 0
$row137[0,24] = 0
$row137[0,25] = 0
$row137[0,26] = 0
$row137[0,27] = $null
$row137[0,28] = $null
$ws.Range("A137:AC137").Value = $row137

$row138 = New-Object "object[,]" 1,29
$row138[0,0] = 136
$row138[0,1] = 6859032
$row138[0,2] = 'Denmark Division 2'
$row138[0,3] = 'Denmark Division 2'
$row138[0,4] = 45383.375
$row138[0,5] = 'Middelfart'
$row138[0,6] = 'Thisted FC'
$row138[0,7] = $null
$row138[0,8] = $null
$row138[0,9] = $null
$row138[0,10] = 1.4
$row138[0,11] = 4.333
$row138[0,12] = 6.5
$row138[0,13] = 1.363
$row138[0,14] = 4.333
$row138[0,15] = 7
$row138[0,16] = -1.25
$row138[0,17] = 1.9
$row138[0,18] = 1.9
$row138[0,19] = 2.75
$row138[0,20] = 1.9
$row138[0,21] = 1.9
$row138[0,22] = 0
$row138[0,23] = 0
$row138[0,24] = 0
$row138[0,25] = 0
$row138[0,26] = 0
$row138[0,27] = $null
$row138[0,28] = $null
$ws.Range("A138:AC138").Value = $row138

$row139 = New-Object "object[,]" 1,29
$row139[0,0] = 137
$row139[0,1] = 6859029
$row139[0,2] = 'Denmark Division 2'
$row139[0,3] = 'Denmark Division 2'
$row139[0,4] = 45383.375
$row139[0,5] = 'Brabrand'
$row139[0,6] = 'Esbjerg'
$row139[0,7] = $null
$row139[0,8] = $null
$row139[0,9] = $null
$row139[0,10] = 6
$row139[0,11] = 4.5
$row139[0,12] = 1.4
$row139[0,13] = 8.5
$row139[0,14] = 5.5
$row139[0,15] = 1.25
$row139[0,16] = 1.75
$row139[0,17] = 1.85
$row139[0,18] = 1.95
$row139[0,19] = 3
$row139[0,20] = 1.825
$row139[0,21] = 1.975
$row139[0,22] = 0
$row139[0,23] = 0
$row139[0,24] = 0
$row139[0,25] = 0
$row139[0,26] = 0
$row139[0,27] = $null
$row139[0,28] = $null
$ws.Range("A139:AC139").Value = $row139

$row140 = New-Object "object[,]" 1,29
$row140[0,0] = 138
$row140[0,1] = 6859025
$row140[0,2] = 'Denmark Division 2'
$row140[0,3] = 'Denmark Division 2'
$row140[0,4] = 45383.58333333334
$row140[0,5] = 'Fremad Amager'
$row140[0,6] = 'AB Copenhagen'
$row140[0,7] = $null
$row140[0,8] = $null
$row140[0,9] = $null
$row140[0,10] = 3
$row140[0,11] = 3.4
$row140[0,12] = 2.1
$row140[0,13] = 3
$row140[0,14] = 3.4
$row140[0,15] = 2.1
$row140[0,16] = 0.25
$row140[0,17] = 1.9
$row140[0,18] = 1.9
$row140[0,19] = 2.5
$row140[0,20] = 1.9
$row140[0,21] = 1.9
$row140[0,22] = 0
$row140[0,23] = 0
$row140[0,24] = 0
$row140[0,25] = 0
$row140[0,26] = 0
$row140[0,27] = $null
$row140[0,28] = $null
$ws.Range("A140:AC140").Value = $row140

Write-Output "done"
